$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A74 with refined timestamp value
$ws.Cells.Item(74, 1).Value = 44387.76718575579

# Add new row 75 with data
$ws.Cells.Item(75, 1).Value = 44388.76741696225
$ws.Cells.Item(75, 2).Value = 79612
$ws.Cells.Item(75, 3).Value = 67242
$ws.Cells.Item(75, 4).Value = 3678
$ws.Cells.Item(75, 5).Value = 2210
$ws.Cells.Item(75, 6).Value = 1579
$ws.Cells.Item(75, 7).Value = 21149
$ws.Cells.Item(75, 8).Value = 1633
$ws.Cells.Item(75, 9).Value = 900
$ws.Cells.Item(75, 10).Value = 208

# Ensure the new row's date cell (A75) uses the same style as A74 (date/time format)
$ws.Cells.Item(75, 1).NumberFormat = $ws.Cells.Item(74, 1).NumberFormat
